# Program Management update: refresh rows 2-13 on the "Data" sheet with the
# new 20220512 (14707) Cocci cartridge batch — Result ID (A), Lab Sample ID
# (E) and Cartridge ID (T) are refreshed for each row; all other columns are
# left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$labSampleId = "20220512-Cocci-14707Updt"
$cartridgeId = "TestCartridge4707"

$resultIds = @{
    2  = "A1425901"
    3  = "A1425902"
    4  = "A1425903"
    5  = "A1425904"
    6  = "A1425905"
    7  = "A1425906"
    8  = "A1425907"
    9  = "A1425908"
    10 = "A1425909"
    11 = "A1425910"
    12 = "A1425911"
    13 = "A1425912"
}

foreach ($row in 2..13) {
    $ws.Range("A$row").Value = $resultIds[$row]
    $ws.Range("E$row").Value = $labSampleId
    $ws.Range("T$row").Value = $cartridgeId
}
